$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in grades for assignments 3 CPP, 4 UE, 5 UE (row 2) ---
# Assignment 3 (CPP) block: N2=Assignment, O2=Grade, P2=Comments
$ws.Range("N2").Value = "3 CPP"
$ws.Range("O2").Value = 75
$ws.Range("P2").Value = "See my comments  under the folder Shahar_Comments"

# Assignment 4 (UE) block: R2=Assignment, S2=Grade, T2=Comments
$ws.Range("T2").Value = "Excellent!"
$ws.Range("R2").Value = "4 UE"
$ws.Range("S2").Value = 100

# Assignment 5 (UE) block: V2=Assignment, W2=Grade, X2=Comments
$ws.Range("V2").Value = "5 UE"
$ws.Range("W2").Value = 33
$ws.Range("X2").Value = "did not create another actor type like Arrow and Target (see items 2 and 3 in class 5 HW)"

# --- Extend the grade table template six more blocks to the right (AC:AZ) ---
# Column block Y:AB (spacer + Assignment/Grade/Comments) is still empty;
# copy it across so the same formatting/template repeats for future assignments.
$src = $ws.Range("Y1:AB13")
$src.Copy($ws.Range("AC1:AF13"))
$src.Copy($ws.Range("AG1:AJ13"))
$src.Copy($ws.Range("AK1:AN13"))
$src.Copy($ws.Range("AO1:AR13"))
$src.Copy($ws.Range("AS1:AV13"))
$src.Copy($ws.Range("AW1:AZ13"))

# --- Widen column X (Comments for assignment 4) so the longer comment is readable ---
$ws.Columns.Item(24).ColumnWidth = 26.28515625

# --- Make header row & row 2 taller to fit the new wrapped comments ---
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 60

# --- Restore the view: scroll right a bit and select P12 ---
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("P12").Select()
